# Weekly update: insert a new price record at the top of the
# "Vega Modelo de Temuco - Espinaca" table (row 143), pushing the
# existing rows 143-152 down to 144-153.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 143; this shifts rows 143:152 down to 144:153
# and keeps the number formatting (date style) consistent with neighboring rows.
$ws.Rows("143:143").Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A143").Value = 10
$ws.Range("B143").Value = "Vega Modelo de Temuco"
$ws.Range("C143").Value = "La Araucanía"
$ws.Range("D143").Value = 44746
$ws.Range("E143").Value = 9
$ws.Range("F143").Value = 100112012
$ws.Range("G143").Value = "Espinaca"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 85
$ws.Range("K143").Value = 10000
$ws.Range("L143").Value = 10000
$ws.Range("M143").Value = 10000
$ws.Range("N143").Value = "$/docena de atados"
$ws.Range("O143").Value = "Región de La Araucanía"
$ws.Range("P143").Value = 3333
$ws.Range("Q143").Value = 3
$ws.Range("R143").Value = "Hortaliza"
